$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: paragraph 1 ("This is a Microsoft word document.") gets two
# trailing spaces appended to the existing run, followed by three new runs
# (split the way the original author split them) containing a red
# (C00000) annotation: "(This is a change - Version for branch alternate)"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1EndMark = $p1.Range.End - 1

$spacesIp = $d.Range($p1EndMark, $p1EndMark)
$spacesIp.InsertAfter("  ")
$afterSpaces = $p1EndMark + 2

$enDash = [char]0x2013

$seg1 = "(This is a change " + $enDash + " Ve"
$ip1 = $d.Range($afterSpaces, $afterSpaces)
$ip1.InsertAfter($seg1)
$seg1End = $afterSpaces + $seg1.Length
$seg1Range = $d.Range($afterSpaces, $seg1End)
$seg1Range.Font.Color = 192

$seg2 = "rsion for branch alternate"
$ip2 = $d.Range($seg1End, $seg1End)
$ip2.InsertAfter($seg2)
$seg2End = $seg1End + $seg2.Length
$seg2Range = $d.Range($seg1End, $seg2End)
$seg2Range.Font.Color = 192

$seg3 = ")"
$ip3 = $d.Range($seg2End, $seg2End)
$ip3.InsertAfter($seg3)
$seg3End = $seg2End + $seg3.Length
$seg3Range = $d.Range($seg2End, $seg3End)
$seg3Range.Font.Color = 192

# ---------------------------------------------------------------------------
# Edit 2: the blank paragraph right after "It will be treated..." (was a
# bare <w:p/>) becomes a shaded (F9F9F9) bold-Calibri heading-style blank
# paragraph (still no run/text, just paragraph mark formatting).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="202122"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($p3xml)

# ---------------------------------------------------------------------------
# Edit 3: the final paragraph ("ank God almighty, we are free at last.",
# styled NormalWeb) is wiped to a completely bare <w:p/> with no pPr/style
# and no runs.
# ---------------------------------------------------------------------------
$countBefore = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($countBefore)
$lastXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$lastPara.Range.InsertXML($lastXml)

# InsertXML at the very end of the story inserts a fresh blank paragraph
# *before* the original (now content-less, but still pPr-carrying) last
# paragraph instead of replacing it in place, so the paragraph count grows
# by one. Collapse the pair back down to a single bare paragraph by
# deleting the old paragraph mark (and its leftover pPr/style) that now
# trails the freshly inserted blank one.
$countAfter = $d.Paragraphs.Count
if ($countAfter -gt $countBefore) {
    $newLastPara = $d.Paragraphs.Item($countAfter)
    $prevPara = $d.Paragraphs.Item($countAfter - 1)
    $seam = $d.Range($prevPara.Range.End - 1, $newLastPara.Range.End)
    $seam.Delete()
}

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
